$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("D").Insert()
$ws.Range("D7").Value = 99999
$ws.Range("D7").NumberFormat = $ws.Range("E7").NumberFormat
Write-Host "done"
